# started processing factor data into modeling data
#
# The "Emerging" quartile block collapses from five columns (quartiles
# 0-4) down to four (quartiles 0-3), matching the "Developed" block.
# Concretely: drop the last quartile column from each block (which
# shifts "Emerging" left to sit right after "Developed" and
# auto-resizes both header merges), then refresh the Forward Return
# sample row with the newly processed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last column of the "Developed" block (F). Excel shifts
# everything to its right one column left (so "Emerging" moves from
# G1 to F1, and the rest of that block follows), and auto-shrinks the
# B1:F1 header merge down to B1:E1.
$ws.Range("F1").EntireColumn.Delete()

# Delete the (now-shifted) last column of the "Emerging" block, which
# landed on J after the first delete. This auto-shrinks the F1:J1
# header merge down to F1:I1, leaving a uniform 4-column block on each
# side, and shrinks the used range from A1:K4 to A1:I4.
$ws.Range("J1").EntireColumn.Delete()

# Refresh the "Forward Return (one month, 21-trading day)" sample row
# with the newly processed values.
$ws.Range("B4").Value2 = 0.00744217190584827
$ws.Range("C4").Value2 = 0.008611320318489409
$ws.Range("D4").Value2 = 0.005006122068839754
$ws.Range("E4").Value2 = 0.008297760554472029
$ws.Range("F4").Value2 = 0.01000803989706718
$ws.Range("G4").Value2 = 0.01080061093467913
$ws.Range("H4").Value2 = 0.01087692005137879
$ws.Range("I4").Value2 = 0.01718382353528924
